$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows per weekly refresh of market data
$ws.Cells.Item(2, 4).Value2 = 44413
$ws.Cells.Item(2, 11).Value2 = 25000
$ws.Cells.Item(2, 13).Value2 = 25000
$ws.Cells.Item(2, 16).Value2 = 2500
$ws.Cells.Item(3, 4).Value2 = 44382
$ws.Cells.Item(3, 10).Value2 = 50
$ws.Cells.Item(3, 11).Value2 = 25000
$ws.Cells.Item(3, 12).Value2 = 25000
$ws.Cells.Item(3, 13).Value2 = 25000
$ws.Cells.Item(3, 16).Value2 = 2500
$ws.Cells.Item(4, 4).Value2 = 44432
$ws.Cells.Item(4, 10).Value2 = 15
$ws.Cells.Item(4, 11).Value2 = 27000
$ws.Cells.Item(4, 12).Value2 = 27000
$ws.Cells.Item(4, 13).Value2 = 27000
$ws.Cells.Item(4, 16).Value2 = 2700
$ws.Cells.Item(5, 4).Value2 = 44428
$ws.Cells.Item(5, 10).Value2 = 30
$ws.Cells.Item(6, 4).Value2 = 44454
$ws.Cells.Item(6, 10).Value2 = 80
$ws.Cells.Item(7, 4).Value2 = 44421
$ws.Cells.Item(7, 10).Value2 = 55
$ws.Cells.Item(8, 4).Value2 = 44379
$ws.Cells.Item(8, 10).Value2 = 35
$ws.Cells.Item(8, 11).Value2 = 22000
$ws.Cells.Item(8, 12).Value2 = 22000
$ws.Cells.Item(8, 13).Value2 = 22000
$ws.Cells.Item(8, 16).Value2 = 2200
$ws.Cells.Item(9, 4).Value2 = 44384
$ws.Cells.Item(9, 10).Value2 = 40
$ws.Cells.Item(10, 4).Value2 = 44349
$ws.Cells.Item(10, 10).Value2 = 45
$ws.Cells.Item(11, 4).Value2 = 44396
$ws.Cells.Item(11, 10).Value2 = 20
$ws.Cells.Item(11, 11).Value2 = 25000
$ws.Cells.Item(11, 12).Value2 = 25000
$ws.Cells.Item(11, 13).Value2 = 25000
$ws.Cells.Item(11, 16).Value2 = 2500
$ws.Cells.Item(12, 4).Value2 = 44412
$ws.Cells.Item(12, 10).Value2 = 50
$ws.Cells.Item(13, 4).Value2 = 44431
$ws.Cells.Item(13, 10).Value2 = 65
$ws.Cells.Item(14, 4).Value2 = 44348
$ws.Cells.Item(14, 10).Value2 = 3
$ws.Cells.Item(15, 4).Value2 = 44435
$ws.Cells.Item(15, 10).Value2 = 185
$ws.Cells.Item(15, 11).Value2 = 25000
$ws.Cells.Item(15, 12).Value2 = 27000
$ws.Cells.Item(15, 13).Value2 = 25162
$ws.Cells.Item(15, 16).Value2 = 2516
$ws.Cells.Item(16, 4).Value2 = 44449
$ws.Cells.Item(16, 10).Value2 = 12
$ws.Cells.Item(17, 4).Value2 = 44355
$ws.Cells.Item(17, 10).Value2 = 25
$ws.Cells.Item(17, 11).Value2 = 23000
$ws.Cells.Item(17, 12).Value2 = 24000
$ws.Cells.Item(17, 13).Value2 = 23400
$ws.Cells.Item(17, 16).Value2 = 2340
$ws.Cells.Item(18, 4).Value2 = 44446
$ws.Cells.Item(18, 10).Value2 = 40
$ws.Cells.Item(18, 11).Value2 = 27000
$ws.Cells.Item(18, 12).Value2 = 27000
$ws.Cells.Item(18, 13).Value2 = 27000
$ws.Cells.Item(18, 16).Value2 = 2700
$ws.Cells.Item(21, 4).Value2 = 44417
$ws.Cells.Item(21, 10).Value2 = 15
$ws.Cells.Item(23, 4).Value2 = 44434
$ws.Cells.Item(23, 10).Value2 = 55
$ws.Cells.Item(23, 11).Value2 = 25000
$ws.Cells.Item(23, 12).Value2 = 25000
$ws.Cells.Item(23, 13).Value2 = 25000
$ws.Cells.Item(23, 16).Value2 = 2500
$ws.Cells.Item(24, 4).Value2 = 44452
$ws.Cells.Item(24, 10).Value2 = 80
$ws.Cells.Item(24, 11).Value2 = 25000
$ws.Cells.Item(24, 12).Value2 = 25000
$ws.Cells.Item(24, 13).Value2 = 25000
$ws.Cells.Item(24, 16).Value2 = 2500
$ws.Cells.Item(25, 4).Value2 = 44441
$ws.Cells.Item(25, 10).Value2 = 70
$ws.Cells.Item(26, 4).Value2 = 44405
$ws.Cells.Item(26, 10).Value2 = 40
$ws.Cells.Item(26, 11).Value2 = 25000
$ws.Cells.Item(26, 12).Value2 = 25000
$ws.Cells.Item(26, 13).Value2 = 25000
$ws.Cells.Item(26, 16).Value2 = 2500
$ws.Cells.Item(27, 4).Value2 = 44426
$ws.Cells.Item(27, 10).Value2 = 30
$ws.Cells.Item(28, 4).Value2 = 44400
$ws.Cells.Item(28, 10).Value2 = 12
$ws.Cells.Item(28, 11).Value2 = 24000
$ws.Cells.Item(28, 12).Value2 = 24000
$ws.Cells.Item(28, 13).Value2 = 24000
$ws.Cells.Item(28, 16).Value2 = 2400
$ws.Cells.Item(29, 4).Value2 = 44392
$ws.Cells.Item(29, 10).Value2 = 25
$ws.Cells.Item(29, 11).Value2 = 24000
$ws.Cells.Item(29, 12).Value2 = 24000
$ws.Cells.Item(29, 13).Value2 = 24000
$ws.Cells.Item(29, 16).Value2 = 2400
$ws.Cells.Item(30, 4).Value2 = 44453
$ws.Cells.Item(30, 10).Value2 = 40
$ws.Cells.Item(30, 11).Value2 = 27000
$ws.Cells.Item(30, 12).Value2 = 27000
$ws.Cells.Item(30, 13).Value2 = 27000
$ws.Cells.Item(30, 16).Value2 = 2700
$ws.Cells.Item(31, 4).Value2 = 44448
$ws.Cells.Item(31, 10).Value2 = 15
$ws.Cells.Item(32, 4).Value2 = 44354
$ws.Cells.Item(32, 10).Value2 = 30
$ws.Cells.Item(32, 11).Value2 = 24000
$ws.Cells.Item(32, 12).Value2 = 24000
$ws.Cells.Item(32, 13).Value2 = 24000
$ws.Cells.Item(32, 16).Value2 = 2400
$ws.Cells.Item(34, 4).Value2 = 44350
$ws.Cells.Item(34, 10).Value2 = 40
$ws.Cells.Item(34, 11).Value2 = 24000
$ws.Cells.Item(34, 13).Value2 = 24375
$ws.Cells.Item(34, 16).Value2 = 2438
$ws.Cells.Item(36, 4).Value2 = 44389
$ws.Cells.Item(36, 10).Value2 = 65
$ws.Cells.Item(37, 4).Value2 = 44447
$ws.Cells.Item(37, 10).Value2 = 30
$ws.Cells.Item(37, 11).Value2 = 27000
$ws.Cells.Item(37, 12).Value2 = 27000
$ws.Cells.Item(37, 13).Value2 = 27000
$ws.Cells.Item(37, 16).Value2 = 2700
$ws.Cells.Item(38, 4).Value2 = 44390
$ws.Cells.Item(38, 11).Value2 = 25000
$ws.Cells.Item(38, 12).Value2 = 25000
$ws.Cells.Item(38, 13).Value2 = 25000
$ws.Cells.Item(38, 16).Value2 = 2500
$ws.Cells.Item(39, 4).Value2 = 44386
$ws.Cells.Item(39, 10).Value2 = 20
$ws.Cells.Item(39, 11).Value2 = 25000
$ws.Cells.Item(39, 12).Value2 = 25000
$ws.Cells.Item(39, 13).Value2 = 25000
$ws.Cells.Item(39, 16).Value2 = 2500
$ws.Cells.Item(40, 4).Value2 = 44433
$ws.Cells.Item(40, 10).Value2 = 25
$ws.Cells.Item(40, 11).Value2 = 25000
$ws.Cells.Item(40, 12).Value2 = 25000
$ws.Cells.Item(40, 13).Value2 = 25000
$ws.Cells.Item(40, 16).Value2 = 2500
$ws.Cells.Item(41, 4).Value2 = 44371
$ws.Cells.Item(41, 10).Value2 = 50
$ws.Cells.Item(42, 4).Value2 = 44397
$ws.Cells.Item(42, 10).Value2 = 30
$ws.Cells.Item(42, 11).Value2 = 27000
$ws.Cells.Item(42, 12).Value2 = 27000
$ws.Cells.Item(42, 13).Value2 = 27000
$ws.Cells.Item(42, 16).Value2 = 2700
$ws.Cells.Item(44, 4).Value2 = 44372
$ws.Cells.Item(44, 10).Value2 = 20
$ws.Cells.Item(45, 4).Value2 = 44356
$ws.Cells.Item(45, 10).Value2 = 15
$ws.Cells.Item(45, 11).Value2 = 24000
$ws.Cells.Item(45, 12).Value2 = 24000
$ws.Cells.Item(45, 13).Value2 = 24000
$ws.Cells.Item(45, 16).Value2 = 2400
$ws.Cells.Item(46, 4).Value2 = 44365
$ws.Cells.Item(46, 10).Value2 = 85
$ws.Cells.Item(46, 11).Value2 = 22000
$ws.Cells.Item(46, 12).Value2 = 22000
$ws.Cells.Item(46, 13).Value2 = 22000
$ws.Cells.Item(46, 16).Value2 = 2200
$ws.Cells.Item(47, 4).Value2 = 44410
$ws.Cells.Item(48, 4).Value2 = 44411
$ws.Cells.Item(48, 10).Value2 = 40
$ws.Cells.Item(49, 4).Value2 = 44376
$ws.Cells.Item(49, 10).Value2 = 45
$ws.Cells.Item(49, 11).Value2 = 23000
$ws.Cells.Item(49, 12).Value2 = 23000
$ws.Cells.Item(49, 13).Value2 = 23000
$ws.Cells.Item(49, 16).Value2 = 2300

# Append new row 50 with latest weekly record
$ws.Cells.Item(50, 1).Value2 = 10
$ws.Cells.Item(50, 2).Value2 = 'Vega Modelo de Temuco'
$ws.Cells.Item(50, 3).Value2 = 'La Araucanía'
$ws.Cells.Item(50, 4).Value2 = 44425
$ws.Cells.Item(50, 5).Value2 = 9
$ws.Cells.Item(50, 6).Value2 = 100112035
$ws.Cells.Item(50, 7).Value2 = 'Bruselas (repollito)'
$ws.Cells.Item(50, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(50, 9).Value2 = 'Primera'
$ws.Cells.Item(50, 10).Value2 = 30
$ws.Cells.Item(50, 11).Value2 = 25000
$ws.Cells.Item(50, 12).Value2 = 25000
$ws.Cells.Item(50, 13).Value2 = 25000
$ws.Cells.Item(50, 14).Value2 = '$/malla 10 kilos'
$ws.Cells.Item(50, 15).Value2 = 'Provincia de Quillota'
$ws.Cells.Item(50, 16).Value2 = 2500
$ws.Cells.Item(50, 17).Value2 = 10
$ws.Cells.Item(50, 18).Value2 = 'Hortaliza'
$ws.Range("D50").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
